$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto data.
# D-column price cells are kept as plain Text (format "@") so values such as
# "1.001", "291.86", "0.3720" etc. are preserved exactly as strings, matching
# the source inline-string cells, instead of being auto-converted to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.488.56"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.573.12"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.86"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3720"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.97"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3399"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.33"
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.053"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.967"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.570.29"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.74"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06764"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("E21").Value = "  +1.99%  "
$ws.Range("E22").Value = "  -1.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.19"
$ws.Range("E23").Value = "  +1.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.482.78"
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.361"
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.638"
$ws.Range("E26").Value = "  -1.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.04"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.42"
$ws.Range("E28").Value = "  +1.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.053"
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.22"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.747.97"
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.085"
$ws.Range("E32").Value = "  +9.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.243"
$ws.Range("E33").Value = "  +3.00%  "
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.778"
$ws.Range("E35").Value = "  -3.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08364"
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2308"
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.339"
$ws.Range("E39").Value = "  -3.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06529"
$ws.Range("E40").Value = "  +0.43%  "
$ws.Range("E41").Value = "  +1.53%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("E43").Value = "  -1.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.08"
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.816"
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.70"
$ws.Range("E48").Value = "  +5.18%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.219"
$ws.Range("E50").Value = "  -4.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07335"
$ws.Range("E51").Value = "  +0.26%  "
